$d = $word.ActiveDocument

# --- Change 1: merge the "Date:" / "12" / "-Feb-2024" runs into one run reading
#     "Date:12-Feb-2024". Finding and replacing with the very same text makes Word
#     coalesce the matched runs into a single run. ---
$d.Content.Find.Execute("Date:12-Feb-2024", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Date:12-Feb-2024", 2)

# --- Change 2: merge the "From hour " / "10" / "am to " / "12:30" runs into one run
#     reading "From hour 10am to 12:30". ---
$d.Content.Find.Execute("From hour 10am to 12:30", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "From hour 10am to 12:30", 2)

# --- Change 3: "Tic tac toe" becomes "Linux project was presented" as its own run,
#     leaving the remaining " mini project was to be presented to the board with
#     program documentation " text as a second, separate run in the same paragraph
#     (still followed by the existing bookmarkEnd). Wrapping the located text in a
#     temporary bookmark before rewriting it keeps the new run from being re-merged
#     with the untouched remainder when the bookmark is removed again. ---
$hit = $d.Content
$hit.Find.Execute("Tic tac toe", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$d.Bookmarks.Add("__tmp_split__", $hit)
$bm = $d.Bookmarks.Item("__tmp_split__")
$bm.Range.Text = "Linux project was presented"
$d.Bookmarks.Item("__tmp_split__").Delete()
